$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2 (Genoa vs Como): odds update ---
$ws.Range("O2").Value = 1.44
$ws.Range("P2").Value = 2.75

# --- Insert two blank rows at position 3, pushing old rows 3-5 down to 5-7 ---
$ws.Rows("3:4").Insert()

# --- New row 3: Argentina - Talleres Cordoba vs Lanus ---
$ws.Range("A3").Value = 'h0hIZWhT'
$ws.Range("B3").NumberFormat = "@"
$ws.Range("B3").Value = '07/11/2024'
$ws.Range("B3").Style = "Normal"
$ws.Range("C3").Value = '18:45'
$ws.Range("D3").Value = 'ARGENTINA - TORNEO BETANO'
$ws.Range("E3").Value = 'Talleres Cordoba'
$ws.Range("F3").Value = 'Lanus'
$ws.Range("G3").Value = 1.9
$ws.Range("H3").Value = 3.25
$ws.Range("I3").Value = 4.5
$ws.Range("J3").Value = 2.63
$ws.Range("K3").Value = 1.95
$ws.Range("L3").Value = 5
$ws.Range("M3").Value = 1.1
$ws.Range("N3").Value = 7
$ws.Range("O3").Value = 1.44
$ws.Range("P3").Value = 2.63
$ws.Range("Q3").Value = 2.5
$ws.Range("R3").Value = 1.5
$ws.Range("S3").Value = 1.53
$ws.Range("T3").Value = 2.38
$ws.Range("U3").Value = 2.2
$ws.Range("V3").Value = 1.62
$ws.Range("W3").Value = 5.5
$ws.Range("X3").Value = 7.5
$ws.Range("Y3").Value = 9.5
$ws.Range("Z3").Value = 15
$ws.Range("AA3").Value = 19
$ws.Range("AB3").Value = 41
$ws.Range("AC3").Value = 6.5
$ws.Range("AD3").Value = 6.5
$ws.Range("AE3").Value = 21
$ws.Range("AF3").Value = 81
$ws.Range("AG3").Value = 501
$ws.Range("AH3").Value = 9.5
$ws.Range("AI3").Value = 21
$ws.Range("AJ3").Value = 17
$ws.Range("AK3").Value = 51
$ws.Range("AL3").Value = 41
$ws.Range("AM3").Value = 51
$ws.Range("AN3").Value = 3.75
$ws.Range("AO3").Value = 11
$ws.Range("AP3").Value = 26
$ws.Range("AQ3").Value = 41
$ws.Range("AR3").Value = 67
$ws.Range("AS3").Value = 251
$ws.Range("AT3").Value = 2.38
$ws.Range("AU3").Value = 9.5
$ws.Range("AV3").Value = 81
$ws.Range("AW3").Value = 6
$ws.Range("AX3").Value = 26
$ws.Range("AY3").Value = 41
$ws.Range("AZ3").Value = 101
$ws.Range("BA3").Value = 151
$ws.Range("BB3").Value = 401
$ws.Range("BC3").Value = 126
$ws.Range("BD3").Value = 126

# --- New row 4: Colombia - Junior vs Millonarios ---
$ws.Range("A4").Value = '0KUiA8fL'
$ws.Range("B4").NumberFormat = "@"
$ws.Range("B4").Value = '07/11/2024'
$ws.Range("B4").Style = "Normal"
$ws.Range("C4").Value = '18:50'
$ws.Range("D4").Value = 'COLOMBIA - PRIMERA A'
$ws.Range("E4").Value = 'Junior'
$ws.Range("F4").Value = 'Millonarios'
$ws.Range("G4").Value = 2.2
$ws.Range("H4").Value = 3
$ws.Range("I4").Value = 3.6
$ws.Range("J4").Value = 3
$ws.Range("K4").Value = 1.95
$ws.Range("L4").Value = 4.33
$ws.Range("M4").Value = 1.1
$ws.Range("N4").Value = 7
$ws.Range("O4").Value = 1.5
$ws.Range("P4").Value = 2.5
$ws.Range("Q4").Value = 2.5
$ws.Range("R4").Value = 1.5
$ws.Range("S4").Value = 1.57
$ws.Range("T4").Value = 2.25
$ws.Range("U4").Value = 2.1
$ws.Range("V4").Value = 1.67
$ws.Range("W4").Value = 6
$ws.Range("X4").Value = 9
$ws.Range("Y4").Value = 10
$ws.Range("Z4").Value = 21
$ws.Range("AA4").Value = 21
$ws.Range("AB4").Value = 41
$ws.Range("AC4").Value = 6.5
$ws.Range("AD4").Value = 6
$ws.Range("AE4").Value = 19
$ws.Range("AF4").Value = 67
$ws.Range("AG4").Value = 201
$ws.Range("AH4").Value = 8.5
$ws.Range("AI4").Value = 17
$ws.Range("AJ4").Value = 13
$ws.Range("AK4").Value = 41
$ws.Range("AL4").Value = 34
$ws.Range("AM4").Value = 41
$ws.Range("AN4").Value = 4
$ws.Range("AO4").Value = 13
$ws.Range("AP4").Value = 29
$ws.Range("AQ4").Value = 41
$ws.Range("AR4").Value = 81
$ws.Range("AS4").Value = 251
$ws.Range("AT4").Value = 2.25
$ws.Range("AU4").Value = 9
$ws.Range("AV4").Value = 67
$ws.Range("AW4").Value = 5.5
$ws.Range("AX4").Value = 21
$ws.Range("AY4").Value = 34
$ws.Range("AZ4").Value = 81
$ws.Range("BA4").Value = 126
$ws.Range("BB4").Value = 351
$ws.Range("BC4").Value = 126
$ws.Range("BD4").Value = 126

# --- Row 5 (Egypt - ZED vs Al Ahly, now shifted down): odds refresh ---
$ws.Range("G5").Value = 6.1
$ws.Range("H5").Value = 3.75
$ws.Range("I5").Value = 1.55
$ws.Range("J5").Value = 5.9
$ws.Range("K5").Value = 2.22
$ws.Range("L5").Value = 2.05
$ws.Range("M5").Value = 1.07
$ws.Range("N5").Value = 7
$ws.Range("O5").Value = 1.32
$ws.Range("P5").Value = 3.1
$ws.Range("Q5").Value = 1.95
$ws.Range("R5").Value = 1.8
$ws.Range("S5").Value = 1.39
$ws.Range("T5").Value = 2.77
$ws.Range("U5").Value = 2
$ws.Range("V5").Value = 1.72
$ws.Range("W5").Value = 14
$ws.Range("X5").Value = 35
$ws.Range("Y5").Value = 19.5
$ws.Range("Z5").Value = 120
$ws.Range("AA5").Value = 70
$ws.Range("AB5").Value = 75
$ws.Range("AC5").Value = 7
$ws.Range("AD5").Value = 7.3
$ws.Range("AE5").Value = 19
$ws.Range("AF5").Value = 100
$ws.Range("AG5").Value = 900
$ws.Range("AH5").Value = 6.1
$ws.Range("AI5").Value = 6.7
$ws.Range("AJ5").Value = 8
$ws.Range("AK5").Value = 10.75
$ws.Range("AL5").Value = 13
$ws.Range("AM5").Value = 30
$ws.Range("AN5").Value = 7.5
$ws.Range("AO5").Value = 37
$ws.Range("AP5").Value = 40
$ws.Range("AQ5").Value = 250
$ws.Range("AR5").Value = 300
$ws.Range("AS5").Value = 500
$ws.Range("AT5").Value = 2.77
$ws.Range("AU5").Value = 8
$ws.Range("AV5").Value = 80
$ws.Range("AW5").Value = 3.3
$ws.Range("AX5").Value = 7.2
$ws.Range("AY5").Value = 17
$ws.Range("AZ5").Value = 22
$ws.Range("BA5").Value = 55
$ws.Range("BB5").Value = 250
$ws.Range("BC5").Value = 51
$ws.Range("BD5").Value = 51
